$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = 0
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = -4
$ws.Range("F9").Value = -10
$ws.Range("F11").Value = -3
